$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-8 from 45221 to 45224
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 45224
}
